$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet name to reflect the new "through" date
$ws.Name = "Through 2022-10-28"

# Update the column header label (shared string) for the current-month column
$ws.Range("B1").Value = "October 2022 (through October 28)"

# Apply the per-neighborhood count updates for 2022-10-28 data (new incidents on 2022-11-05)
$ws.Range("AF2").Value = 8
$ws.Range("AZ3").Value = 3
$ws.Range("BJ4").Value = 4
$ws.Range("AZ7").Value = 5
$ws.Range("BJ8").Value = 1
$ws.Range("B10").Value = 6
$ws.Range("AP12").Value = 1
$ws.Range("B15").Value = 2
$ws.Range("AF15").Value = 1
$ws.Range("V16").Value = 6
$ws.Range("AP17").Value = 1
$ws.Range("AF18").Value = 2
$ws.Range("AZ18").Value = 5
$ws.Range("BJ18").Value = 4
$ws.Range("BT20").Value = 1
$ws.Range("L24").Value = 6
$ws.Range("AF25").Value = 2
$ws.Range("AZ31").Value = 1
$ws.Range("L42").Value = 3
$ws.Range("V44").Value = 2
$ws.Range("AZ70").Value = 1
$ws.Range("V79").Value = 2
$ws.Range("AZ79").Value = 1
$ws.Range("B81").Value = 1
$ws.Range("AZ95").Value = 2
$ws.Range("AZ97").Value = 2
$ws.Range("B98").Value = 2
